# Insert a new data row before existing row 720 (shifts rows 720-778 down to 721-779)
# so the table grows from A1:R778 to A1:R779, matching a new weekly price entry
# for "Espinaca" at "Mercado Mayorista Lo Valledor de Santiago".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(720).Insert()

$ws.Range("A720").Value = 6
$ws.Range("B720").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C720").Value = "Metropolitana"
$ws.Range("D720").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D720").Value = 45106
$ws.Range("E720").Value = 13
$ws.Range("F720").Value = 100112012
$ws.Range("G720").Value = "Espinaca"
$ws.Range("H720").Value = "Sin especificar"
$ws.Range("I720").Value = "Primera"
$ws.Range("J720").Value = 510
$ws.Range("K720").Value = 6000
$ws.Range("L720").Value = 6500
$ws.Range("M720").Value = 6225
$ws.Range("N720").Value = "$/cuna 10 kilos"
$ws.Range("O720").Value = "Región Metropolitana"
$ws.Range("P720").Value = 622
$ws.Range("Q720").Value = 10
$ws.Range("R720").Value = "Hortaliza"
